# Put the git metadata attributes in the right order: Url, Branch,
# Revision (previously Branch, Revision, Url), add a new blank row
# below them, and move the active tab/selection to the
# "Data repo metadata" sheet (was "Model1s").

$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item("Data repo metadata")

# --- capture current label/value pairs before shuffling them ---
$branchLabel   = $metaSheet.Range("A1").Value2
$branchValue   = $metaSheet.Range("B1").Value2
$revisionLabel = $metaSheet.Range("A2").Value2
$revisionValue = $metaSheet.Range("B2").Value2
$urlLabel      = $metaSheet.Range("A3").Value2
$urlValue      = $metaSheet.Range("B3").Value2

# --- capture the existing per-row data validations (Branch/B1,
#     Revision/B2, Url/B3) so we can re-create them against the rows
#     the fields move to ---
$branchVal = $metaSheet.Range("B1").Validation
$branchValErrorTitle = $branchVal.ErrorTitle
$branchValErrorMessage = $branchVal.ErrorMessage
$branchValInputTitle = $branchVal.InputTitle
$branchValInputMessage = $branchVal.InputMessage

$revisionVal = $metaSheet.Range("B2").Validation
$revisionValErrorTitle = $revisionVal.ErrorTitle
$revisionValErrorMessage = $revisionVal.ErrorMessage
$revisionValInputTitle = $revisionVal.InputTitle
$revisionValInputMessage = $revisionVal.InputMessage

$urlVal = $metaSheet.Range("B3").Validation
$urlValErrorTitle = $urlVal.ErrorTitle
$urlValErrorMessage = $urlVal.ErrorMessage
$urlValInputTitle = $urlVal.InputTitle
$urlValInputMessage = $urlVal.InputMessage

# --- reorder the rows: Url, Branch, Revision ---
$metaSheet.Range("A1").Value2 = $urlLabel
$metaSheet.Range("B1").Value2 = $urlValue
$metaSheet.Range("A2").Value2 = $branchLabel
$metaSheet.Range("B2").Value2 = $branchValue
$metaSheet.Range("A3").Value2 = $revisionLabel
$metaSheet.Range("B3").Value2 = $revisionValue

# --- add a new, blank 4th row to the sheet's used range (matches the
#     row Excel appends below the last data-validated field) ---
$metaSheet.Range("A4:B4").Style = "Normal"
$metaSheet.Rows.Item(4).RowHeight = 15

# --- drop the old validations and recreate them on the rows the
#     fields now occupy: Url -> B1, Branch -> B2, Revision -> B3 ---
$metaSheet.Range("B1").Validation.Delete()
$metaSheet.Range("B2").Validation.Delete()
$metaSheet.Range("B3").Validation.Delete()

$metaSheet.Range("B1").Validation.Add(6, 2, 8, 255)
$v = $metaSheet.Range("B1").Validation
$v.ErrorTitle = $urlValErrorTitle
$v.ErrorMessage = $urlValErrorMessage
$v.InputTitle = $urlValInputTitle
$v.InputMessage = $urlValInputMessage

$metaSheet.Range("B2").Validation.Add(6, 2, 8, 255)
$v = $metaSheet.Range("B2").Validation
$v.ErrorTitle = $branchValErrorTitle
$v.ErrorMessage = $branchValErrorMessage
$v.InputTitle = $branchValInputTitle
$v.InputMessage = $branchValInputMessage

$metaSheet.Range("B3").Validation.Add(6, 2, 8, 255)
$v = $metaSheet.Range("B3").Validation
$v.ErrorTitle = $revisionValErrorTitle
$v.ErrorMessage = $revisionValErrorMessage
$v.InputTitle = $revisionValInputTitle
$v.InputMessage = $revisionValInputMessage

# --- make "Data repo metadata" the active sheet/tab, with the new
#     blank row selected ---
$metaSheet.Activate()
$metaSheet.Range("A4:XFD4").Select()
